$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27; existing rows 27.. shift down to 28..
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly price record.
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value = 44495
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 100112036
$ws.Cells.Item(27, 7).Value = "Caigua"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 130
$ws.Cells.Item(27, 11).Value = 5000
$ws.Cells.Item(27, 12).Value = 6000
$ws.Cells.Item(27, 13).Value = 5500
$ws.Cells.Item(27, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 275
$ws.Cells.Item(27, 17).Value = 20
$ws.Cells.Item(27, 18).Value = "Hortaliza"
